# Bot_feed.xlsx restructuring:
# The single wide header row (A1:J1) is converted into a two-column layout:
#   - A1/B1 stay as-is (the greeting + intro question)
#   - The remaining values that used to live in C1:J1 move down into
#     column B, one per row, starting at B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing values (and the source range address) before we
# start moving anything around.
$sourceCells = @("C1", "D1", "E1", "F1", "G1", "H1", "I1", "J1")
$values = @()
foreach ($addr in $sourceCells) {
    $values += ,($ws.Range($addr).Value())
}

# Remove the old wide cells (clears both value and formatting).
$ws.Range("C1:J1").Clear()

# Write the captured values down column B, starting on row 2, and make
# sure they pick up the same wrapped-text style the header row used.
$destRow = 2
foreach ($val in $values) {
    $cell = $ws.Cells.Item($destRow, 2)
    $cell.Value = $val
    $cell.WrapText = $true
    $destRow = $destRow + 1
}

# Row heights to match the new, taller, per-item rows. Rows 2, 4 and 7
# keep the sheet's default 14.5 height, so they're left untouched.
$ws.Rows("1").RowHeight = 377
$ws.Rows("3").RowHeight = 29
$ws.Rows("5").RowHeight = 304.5
$ws.Rows("6").RowHeight = 87
$ws.Rows("8").RowHeight = 43.5
$ws.Rows("9").RowHeight = 72.5

# Restore the selection Excel leaves behind after this kind of
# cut/paste-down rearrangement.
$ws.Range("B1:J1").Select() | Out-Null
